$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '29.166.76'
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  +0.09%  '
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.833.80'
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -0.22%  '
$c.ClearFormats()

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '0.9994'
$c.ClearFormats()
$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '241.27'
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  +0.58%  '
$c.ClearFormats()

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.6652'
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -2.40%  '
$c.ClearFormats()

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c.ClearFormats()

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.07423'
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  -0.40%  '
$c.ClearFormats()

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.2938'
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -1.94%  '
$c.ClearFormats()

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '22.64'
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -2.36%  '
$c.ClearFormats()

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07738'
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  +1.31%  '
$c.ClearFormats()

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.810.74'
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  -1.57%  '
$c.ClearFormats()

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.986'
$c.ClearFormats()
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -1.11%  '
$c.ClearFormats()

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '0.6691'
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -1.60%  '
$c.ClearFormats()

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '83.00'
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  -5.51%  '
$c.ClearFormats()

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '6.099'
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -0.28%  '
$c.ClearFormats()

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '0.000008400'
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  +2.61%  '
$c.ClearFormats()

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '29.100.20'
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  -0.17%  '
$c.ClearFormats()

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '227.19'
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -1.50%  '
$c.ClearFormats()

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '12.47'
$c.ClearFormats()
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  -0.39%  '
$c.ClearFormats()

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  +0.09%  '
$c.ClearFormats()

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '7.183'
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -2.29%  '
$c.ClearFormats()

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -0.05%  '
$c.ClearFormats()

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '159.61'
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -0.56%  '
$c.ClearFormats()

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '0.1412'
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  -1.73%  '
$c.ClearFormats()

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '8.629'
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -0.74%  '
$c.ClearFormats()

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -0.62%  '
$c.ClearFormats()

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '1.513'
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +0.79%  '
$c.ClearFormats()

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '4.116'
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -3.47%  '
$c.ClearFormats()

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '4.046'
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -2.24%  '
$c.ClearFormats()

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.195'
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  +0.25%  '
$c.ClearFormats()

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '0.05330'
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  -0.58%  '
$c.ClearFormats()

$c = $ws.Range("B33")
$c.NumberFormat = "@"
$c.Value = 'LidoDAOToken'
$c.ClearFormats()
$c = $ws.Range("C33")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c.ClearFormats()
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '1.879'
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  +1.46%  '
$c.ClearFormats()

$c = $ws.Range("B34")
$c.NumberFormat = "@"
$c.Value = 'ImmutableX'
$c.ClearFormats()
$c = $ws.Range("C34")
$c.NumberFormat = "@"
$c.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c.ClearFormats()
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.7578'
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  +0.45%  '
$c.ClearFormats()

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.137'
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  +0.37%  '
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '2.672'
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -0.63%  '
$c.ClearFormats()

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.272.35'
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -2.96%  '
$c.ClearFormats()

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.01799'
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -1.65%  '
$c.ClearFormats()

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '2.735'
$c.ClearFormats()
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  +0.48%  '
$c.ClearFormats()

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.9299'
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -1.70%  '
$c.ClearFormats()

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.08887'
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  +16.12%  '
$c.ClearFormats()

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '5.978'
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -0.41%  '
$c.ClearFormats()

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '1.005'
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  +0.56%  '
$c.ClearFormats()

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '102.84'
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -1.67%  '
$c.ClearFormats()

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '1.956.66'
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  -1.69%  '
$c.ClearFormats()

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.5165'
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -0.42%  '
$c.ClearFormats()

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '1.773'
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  +0.22%  '
$c.ClearFormats()

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -1.85%  '
$c.ClearFormats()

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '63.42'
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -1.42%  '
$c.ClearFormats()

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.05916'
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  -0.50%  '
$c.ClearFormats()

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '6.790'
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -1.24%  '
$c.ClearFormats()
